# Update the LR-pairs data with new TPM-derived values.
# The "Target cluster" pattern (ECs, FAPs, Inflammatory-Mac, MuSCs, Resolving-Mac)
# stays the same for each block of 5 rows, but the "Sending cluster" for the
# first block changes from ECs -> Inflammatory-Mac, and for the second block
# from Inflammatory-Mac -> Resolving-Mac. All numeric measurement columns
# (E through T) are refreshed with newly computed TPM-based values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @{Row=2;  A="Inflammatory-Mac"; B="Ifnb1"; C="Ifnar2"; D="ECs";              E=1; F=0.3333333333333333; G=0.1044623333333333; H=0.313387;   I=0.3334131260000681; J=0.3334131260000681; K=3; L=1; M=27.85292233333333; N=83.55876699999999;  O=0.1175699887262562; P=0.1175699887262562; Q=2.909581257092111;  R=26.186231313829;   S=0.03919937746501384; T=0.03919937746501384}
    @{Row=3;  A="Inflammatory-Mac"; B="Ifnb1"; C="Ifnar2"; D="FAPs";             E=1; F=0.3333333333333333; G=0.1044623333333333; H=0.313387;   I=0.3334131260000681; J=0.3334131260000681; K=3; L=1; M=26.66082666666667; N=79.98248000000001;  O=0.1125380329259528; P=0.1125380329259528; Q=2.785052162195556;  R=25.06546945976001; S=0.03752165735174049; T=0.0375216573517405}
    @{Row=4;  A="Inflammatory-Mac"; B="Ifnb1"; C="Ifnar2"; D="Inflammatory-Mac"; E=1; F=0.3333333333333333; G=0.1044623333333333; H=0.313387;   I=0.3334131260000681; J=0.3334131260000681; K=3; L=1; M=105.665011;         N=316.995033;         O=0.4460226472237104; P=0.4460226472237104; Q=11.03801360075233;  R=99.342122406771;   S=0.1487098050776829;  T=0.1487098050776829}
    @{Row=5;  A="Inflammatory-Mac"; B="Ifnb1"; C="Ifnar2"; D="MuSCs";            E=1; F=0.3333333333333333; G=0.1044623333333333; H=0.313387;   I=0.3334131260000681; J=0.3334131260000681; K=3; L=1; M=4.883238666666666;  N=14.649716;          O=0.02061264193813266; P=0.02061264193813266; Q=0.5101145053435555; R=4.591030548092;    S=0.006872525383712911; T=0.006872525383712913}
    @{Row=6;  A="Inflammatory-Mac"; B="Ifnb1"; C="Ifnar2"; D="Resolving-Mac";    E=1; F=0.3333333333333333; G=0.1044623333333333; H=0.313387;   I=0.3334131260000681; J=0.3334131260000681; K=3; L=1; M=71.84303666666666;  N=215.52911;          O=0.3032566891859479; P=0.3032566891859479; Q=7.504891243952222;  R=67.54402119557001; S=0.1011097607219179;  T=0.1011097607219179}
    @{Row=7;  A="Resolving-Mac";     B="Ifnb1"; C="Ifnar2"; D="ECs";              E=1; F=0.3333333333333333; G=0.2088496666666667; H=0.626549;   I=0.6665868739999319; J=0.6665868739999319; K=3; L=1; M=27.85292233333333; N=83.55876699999999;  O=0.1175699887262562; P=0.1175699887262562; Q=5.817073545009221;  R=52.35366190508299; S=0.07837061126124235; T=0.07837061126124235}
    @{Row=8;  A="Resolving-Mac";     B="Ifnb1"; C="Ifnar2"; D="FAPs";             E=1; F=0.3333333333333333; G=0.2088496666666667; H=0.626549;   I=0.6665868739999319; J=0.6665868739999319; K=3; L=1; M=26.66082666666667; N=79.98248000000001;  O=0.1125380329259528; P=0.1125380329259528; Q=5.568104762391112;  R=50.11294286152001; S=0.07501637557421226; T=0.07501637557421226}
    @{Row=9;  A="Resolving-Mac";     B="Ifnb1"; C="Ifnar2"; D="Inflammatory-Mac"; E=1; F=0.3333333333333333; G=0.2088496666666667; H=0.626549;   I=0.6665868739999319; J=0.6665868739999319; K=3; L=1; M=105.665011;         N=316.995033;         O=0.4460226472237104; P=0.4460226472237104; Q=22.06810232567967;  R=198.612920931117;  S=0.2973128421460275;  T=0.2973128421460275}
    @{Row=10; A="Resolving-Mac";     B="Ifnb1"; C="Ifnar2"; D="MuSCs";            E=1; F=0.3333333333333333; G=0.2088496666666667; H=0.626549;   I=0.6665868739999319; J=0.6665868739999319; K=3; L=1; M=4.883238666666666;  N=14.649716;          O=0.02061264193813266; P=0.02061264193813266; Q=1.019862767787111;  R=9.178764910084;    S=0.01374011655441975; T=0.01374011655441975}
    @{Row=11; A="Resolving-Mac";     B="Ifnb1"; C="Ifnar2"; D="Resolving-Mac";    E=1; F=0.3333333333333333; G=0.2088496666666667; H=0.626549;   I=0.6665868739999319; J=0.6665868739999319; K=3; L=1; M=71.84303666666666;  N=215.52911;          O=0.3032566891859479; P=0.3032566891859479; Q=15.00439426015445;  R=135.03954834139;   S=0.20214692846403;    T=0.20214692846403}
)

foreach ($r in $data) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
